$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting from an existing date row so the new row reuses
# the same style (rather than Excel auto-creating a new number format).
$ws.Range("C6:D6").Copy()
$ws.Range("C7:D7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Add a new row for "Make Your Bed"
$ws.Range("A7").Value = "Make Your Bed"
$ws.Range("B7").Value = "William H. McRaven"
$ws.Range("C7").Value = "1/11/2020"
$ws.Range("D7").Value = "1/12/2020"
$ws.Range("E7").Value = "inspiration;success;self-improvement;motivation"
$ws.Range("F7").Value = "Hard Copy"
$ws.Range("G7").Value = "125 Pages"

# Update the selected cell to reflect the new end-of-data position
$ws.Range("A8").Select()
